# Update the MSME indicator figures on the "Summary" sheet to their more
# precise (two-decimal) values. The source cells hold numeric-looking text
# (shared strings), so a leading apostrophe is used to keep Excel from
# reinterpreting the input as a genuine number -- this preserves the
# original "number stored as text" cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people) -- row 11
$ws.Range("B11").Value = "'66.59"
$ws.Range("C11").Value = "'5.39"
$ws.Range("D11").Value = "'71.98"

# Employment (% of total) -- row 12
$ws.Range("B12").Value = "'41.14"
$ws.Range("C12").Value = "'44.92"
$ws.Range("D12").Value = "'86.06"

# Enterprises (% of total) -- row 14
$ws.Range("B14").Value = "'92.43"
$ws.Range("C14").Value = "'7.48"
$ws.Range("D14").Value = "'99.91"
